$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 2026-01-12 / Hall / Regular (date unchanged)
$ws.Range("B2").Value = "Hall"
$ws.Range("C2").Value = 10.5
$ws.Range("E2").Value = 65
$ws.Range("F2").Value = 682.5

# Row 3 - 2026-01-13 / Patton / Regular (date unchanged)
$ws.Range("B3").Value = "Patton"
$ws.Range("C3").Value = 8.5
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 552.5

# Row 4 - 2026-01-14 / Bahin / Regular (date unchanged)
$ws.Range("B4").Value = "Bahin"
$ws.Range("C4").Value = 10.5
$ws.Range("E4").Value = 65
$ws.Range("F4").Value = 682.5

# Row 5 - 2026-01-15 / Lynn / Regular (date unchanged)
$ws.Range("B5").Value = "Lynn"
$ws.Range("C5").Value = 10.5
$ws.Range("E5").Value = 65
$ws.Range("F5").Value = 682.5

# Row 6 - date changes from 2026-01-16 to 2026-01-15 / Lynn / OT
# Use a helper cell with a text formula + paste-special values so the
# date-like literal lands as plain text (matches the other date cells)
# instead of being auto-converted to a date serial number.
$ws.Range("Z1").Formula = "=""2026-01-15"""
$ws.Range("Z1").Copy()
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("Z1").Clear()
$ws.Range("B6").Value = "Lynn"
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = "OT"
$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 48.75

# Row 7 - 2026-01-16 / Lucas / OT (date unchanged)
$ws.Range("B7").Value = "Lucas"
$ws.Range("C7").Value = 9
$ws.Range("E7").Value = 65
$ws.Range("F7").Value = 877.5

# Row 9 - SUBTOTAL
$ws.Range("C9").Value = 49.5
$ws.Range("D9").Value = "Reg: 40 / OT: 9.5"
$ws.Range("F9").Value = 3526.25
